$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert() | Out-Null
$ws.Columns("N:N").ColumnWidth = $mWidth
$ws.Range("R8").Select() | Out-Null
Write-Output "done"
